$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("F2").Value = 25.99000000000062
$ws.Range("H2").Value = 0.03712376461636779
$ws.Range("I2").Value = 0.03712376461636779
$ws.Range("L2").Value = 7.778967120477413
$ws.Range("M2").Value = '[0.23188970232397388, 15.326044538630851]'
$ws.Range("N2").Value = 0.04363770202042239
$ws.Range("O2").Value = 0.04363770202042239
$ws.Range("P2").Value = -1.761052938949233
$ws.Range("Q2").Value = '[-3.107000542289004, -0.4151053356094625]'
$ws.Range("R2").Value = 0.01149101867804903
$ws.Range("S2").Value = 0.01149101867804903
$ws.Range("T2").Value = 12.69971272682398
$ws.Range("U2").Value = '[8.423061006668952, 16.976364446979]'
$ws.Range("V2").Value = 0.0000003332180111836891
$ws.Range("W2").Value = 0.0000003332180111836891
$ws.Range("X2").Value = 7.284484484484658
$ws.Range("Y2").Value = 1.717057057057098
$ws.Range("Z2").Value = 12.85191191191222

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("F3").Value = 25.99000000000062
$ws.Range("H3").Value = 0.044993535262566
$ws.Range("I3").Value = 0.044993535262566
$ws.Range("L3").Value = 7.566832182600196
$ws.Range("M3").Value = '[0.37733847425350575, 14.756325890946886]'
$ws.Range("N3").Value = 0.0395706370820923
$ws.Range("O3").Value = 0.0395706370820923
$ws.Range("P3").Value = -1.534631846798618
$ws.Range("Q3").Value = '[-2.792526803190927, -0.2767368904063088]'
$ws.Range("R3").Value = 0.01792158802515109
$ws.Range("S3").Value = 0.01792158802515109
$ws.Range("T3").Value = 11.97978420797114
$ws.Range("U3").Value = '[7.631590761899311, 16.32797765404296]'
$ws.Range("V3").Value = 0.000001450744481878985
$ws.Range("W3").Value = 0.000001450744481878985
$ws.Range("X3").Value = 6.347907907908063
$ws.Range("Y3").Value = 1.144704704704736
$ws.Range("Z3").Value = 11.55111111111139

# Row 4
$ws.Range("F4").Value = 25.99000000000062
$ws.Range("H4").Value = 0.2058196928931515
$ws.Range("I4").Value = 0.2058196928931515
$ws.Range("L4").Value = 5.413899366762892
$ws.Range("M4").Value = '[-1.8727430376676892, 12.700541771193473]'
$ws.Range("N4").Value = 0.1415157410291639
$ws.Range("O4").Value = 0.1415157410291639
$ws.Range("P4").Value = 3.050395269251351
$ws.Range("Q4").Value = '[0.018868424345884982, 6.081922114156817]'
$ws.Range("R4").Value = 0.04864983252913135
$ws.Range("S4").Value = 0.04864983252913135
$ws.Range("T4").Value = 11.46431771223958
$ws.Range("U4").Value = '[7.326541862914507, 15.602093561564661]'
$ws.Range("V4").Value = 0.000001304665411261396
$ws.Range("W4").Value = 0.000001304665411261396
$ws.Range("X4").Value = 13.37223223223255
$ws.Range("Y4").Value = 0.8325125125125297
$ws.Range("Z4").Value = 25.91195195195257

# Row 5
$ws.Range("F5").Value = 25.99000000000062
$ws.Range("H5").Value = 0.3713251281015377
$ws.Range("I5").Value = 0.3713251281015377
$ws.Range("L5").Value = 4.277322804433428
$ws.Range("M5").Value = '[-3.973516193680373, 12.52816180254723]'
$ws.Range("N5").Value = 0.3020008087496457
$ws.Range("O5").Value = 0.3020008087496457
$ws.Range("P5").Value = -2.138421425866927
$ws.Range("Q5").Value = '[-5.220264069028085, 0.9434212172942313]'
$ws.Range("R5").Value = 0.169103310259761
$ws.Range("S5").Value = 0.169103310259761
$ws.Range("T5").Value = 11.08486832258084
$ws.Range("U5").Value = '[6.82172856470093, 15.348008080460747]'
$ws.Range("V5").Value = 0.000004162458135947844
$ws.Range("W5").Value = 0.000004162458135947844
$ws.Range("X5").Value = 8.845445445445659
$ws.Range("Y5").Value = -3.902402402402494
$ws.Range("Z5").Value = 21.59329329329381

# Row 6
$ws.Range("F6").Value = 25.99000000000062
$ws.Range("H6").Value = 0.04055194411036545
$ws.Range("I6").Value = 0.04055194411036545
$ws.Range("L6").Value = 8.368638530131163
$ws.Range("M6").Value = '[-0.7734992947367658, 17.51077635499909]'
$ws.Range("N6").Value = 0.07181866941265302
$ws.Range("O6").Value = 0.07181866941265302
$ws.Range("P6").Value = 2.358553043235581
$ws.Range("Q6").Value = '[-0.6037895790683088, 5.320895665539471]'
$ws.Range("R6").Value = 0.115802153556823
$ws.Range("S6").Value = 0.115802153556823
$ws.Range("T6").Value = 12.99385013135113
$ws.Range("U6").Value = '[8.32502541252152, 17.66267485018074]'
$ws.Range("V6").Value = 0.000001198086480469485
$ws.Range("W6").Value = 0.000001198086480469485
$ws.Range("X6").Value = 16.23399399399438
$ws.Range("Y6").Value = 3.98045045045054
$ws.Range("Z6").Value = 28.48753753753822

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 23.55000000000024
$ws.Range("H7").Value = 0.04446429218528991
$ws.Range("I7").Value = 0.04446429218528991
$ws.Range("L7").Value = 7.856797919509274
$ws.Range("M7").Value = '[0.7636205515515488, 14.949975287467]'
$ws.Range("N7").Value = 0.03071401753043679
$ws.Range("O7").Value = 0.03071401753043679
$ws.Range("P7").Value = 1.490605523324887
$ws.Range("Q7").Value = '[0.19497371824080822, 2.7862373284089657]'
$ws.Range("R7").Value = 0.02510039802051223
$ws.Range("S7").Value = 0.02510039802051223
$ws.Range("T7").Value = 12.71669815839536
$ws.Range("U7").Value = '[8.42618663271985, 17.007209684070872]'
$ws.Range("V7").Value = 0.0000003463919420187977
$ws.Range("W7").Value = 0.0000003463919420187977
$ws.Range("X7").Value = 17.96306306306325
$ws.Range("Y7").Value = 13.10690690690704
$ws.Range("Z7").Value = 22.81921921921946

# Row 8
$ws.Range("F8").Value = 23.55000000000024
$ws.Range("H8").Value = 0.7681409970024018
$ws.Range("I8").Value = 0.7681409970024018
$ws.Range("L8").Value = 2.312941079217987
$ws.Range("M8").Value = '[-6.363647563170915, 10.989529721606889]'
$ws.Range("N8").Value = 0.5939781898930474
$ws.Range("O8").Value = 0.5939781898930474
$ws.Range("P8").Value = 1.000026490331886
$ws.Range("Q8").Value = '[-2.1384214258669267, 4.138474406530699]'
$ws.Range("R8").Value = 0.5242793640232137
$ws.Range("S8").Value = 0.5242793640232137
$ws.Range("T8").Value = 12.87231960794077
$ws.Range("U8").Value = '[8.37369431078687, 17.370944905094678]'
$ws.Range("V8").Value = 0.0000007007503237677071
$ws.Range("W8").Value = 0.0000007007503237677071
$ws.Range("X8").Value = 19.801801801802
$ws.Range("Y8").Value = 8.038588588588668
$ws.Range("Z8").Value = 31.56501501501534

# Row 9
$ws.Range("F9").Value = 23.55000000000024
$ws.Range("H9").Value = 0.3485919317534392
$ws.Range("I9").Value = 0.3485919317534392
$ws.Range("L9").Value = 4.578987231014357
$ws.Range("M9").Value = '[-3.8555931284430356, 13.013567590471748]'
$ws.Range("N9").Value = 0.2800275023997014
$ws.Range("O9").Value = 0.2800275023997014
$ws.Range("P9").Value = 2.018921405009657
$ws.Range("Q9").Value = '[-1.1195265111891555, 5.157369321208469]'
$ws.Range("R9").Value = 0.2017044520300906
$ws.Range("S9").Value = 0.2017044520300906
$ws.Range("T9").Value = 12.72454585025878
$ws.Range("U9").Value = '[8.312381389874922, 17.136710310642634]'
$ws.Range("V9").Value = 0.0000006001385852716368
$ws.Range("W9").Value = 0.0000006001385852716368
$ws.Range("X9").Value = 15.98288288288305
$ws.Range("Y9").Value = 4.219669669669713
$ws.Range("Z9").Value = 27.74609609609638

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("F10").Value = 23.55000000000024
$ws.Range("H10").Value = 0.02512234308597927
$ws.Range("I10").Value = 0.02512234308597927
$ws.Range("L10").Value = 7.965629564539131
$ws.Range("M10").Value = '[0.658504112449755, 15.272755016628507]'
$ws.Range("N10").Value = 0.03331727353495983
$ws.Range("O10").Value = 0.03331727353495983
$ws.Range("P10").Value = 1.817658211986887
$ws.Range("Q10").Value = '[0.5723422051585008, 3.062974218815274]'
$ws.Range("R10").Value = 0.005169043508775673
$ws.Range("S10").Value = 0.005169043508775673
$ws.Range("T10").Value = 10.78575910353891
$ws.Range("U10").Value = '[6.812249470368755, 14.759268736709068]'
$ws.Range("V10").Value = 0.000001915240951610997
$ws.Range("W10").Value = 0.000001915240951610997
$ws.Range("X10").Value = 16.73723723723741
$ws.Range("Y10").Value = 12.06966966966979
$ws.Range("Z10").Value = 21.40480480480502

# Row 11
$ws.Range("B11").Value = 0
$ws.Range("F11").Value = 23.55000000000024
$ws.Range("H11").Value = 0.1257091024971805
$ws.Range("I11").Value = 0.1257091024971805
$ws.Range("L11").Value = 6.407620162834623
$ws.Range("M11").Value = '[-1.7846992983654557, 14.5999396240347]'
$ws.Range("N11").Value = 0.1221853242373259
$ws.Range("O11").Value = 0.1221853242373259
$ws.Range("P11").Value = 1.968605606753965
$ws.Range("Q11").Value = '[-0.9811580659860013, 4.91836927949393]'
$ws.Range("R11").Value = 0.185632023898666
$ws.Range("S11").Value = 0.185632023898666
$ws.Range("T11").Value = 12.17077887949602
$ws.Range("U11").Value = '[7.823375550851342, 16.518182208140704]'
$ws.Range("V11").Value = 0.00000107059732701309
$ws.Range("W11").Value = 0.00000107059732701309
$ws.Range("X11").Value = 16.17147147147164
$ws.Range("Y11").Value = 5.115465465465517
$ws.Range("Z11").Value = 27.22747747747776

# Row 12
$ws.Range("F12").Value = 23.55000000000024
$ws.Range("H12").Value = 0.3971275301278157
$ws.Range("I12").Value = 0.3971275301278157
$ws.Range("L12").Value = 4.324393697636097
$ws.Range("M12").Value = '[-4.405201135744912, 13.053988531017106]'
$ws.Range("N12").Value = 0.3237454853844113
$ws.Range("O12").Value = 0.3237454853844113
$ws.Range("P12").Value = 2.232763547596349
$ws.Range("Q12").Value = '[-0.9056843686024632, 5.371211463795161]'
$ws.Range("R12").Value = 0.1588034606640034
$ws.Range("S12").Value = 0.1588034606640034
$ws.Range("T12").Value = 9.799262733204607
$ws.Range("U12").Value = '[5.325861542419984, 14.27266392398923]'
$ws.Range("V12").Value = 0.00006324294189830049
$ws.Range("W12").Value = 0.00006324294189830049
$ws.Range("X12").Value = 15.18138138138154
$ws.Range("Y12").Value = 3.418168168168206
$ws.Range("Z12").Value = 26.94459459459488
